$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 0  # H17: was 2238.5
$ws.Cells.Item(17, 9).Value = 0  # I17: was 2500
$ws.Cells.Item(17, 10).Value = 0  # J17: was 2151.3333
$ws.Cells.Item(17, 11).Value = 0  # K17: was 7500
$ws.Cells.Item(17, 12).Value = 0  # L17: was 6453.999899999999
$ws.Cells.Item(17, 13).ClearContents()  # M17: was -7332
$ws.Cells.Item(17, 14).ClearContents()  # N17: was -6789.999899999999

$ws.Cells.Item(32, 8).Value = 2372.7856  # H32: was 2479.4614
$ws.Cells.Item(32, 9).Value = 1740.75  # I32: was 1992.3334
$ws.Cells.Item(32, 11).Value = 1740.75  # K32: was 1992.3334
$ws.Cells.Item(32, 13).Value = -1414.75  # M32: was -1666.3334

$ws.Cells.Item(33, 8).Value = 12979.75  # H33: was 14699
$ws.Cells.Item(33, 10).Value = 956.6667  # J33: was 962.5
$ws.Cells.Item(33, 12).Value = 956.6667  # L33: was 962.5
$ws.Cells.Item(33, 14).Value = -1414.6667  # N33: was -1420.5

$ws.Cells.Item(40, 8).Value = 4169.2856  # H40: was 4382.846
$ws.Cells.Item(40, 9).Value = 2261.6667  # I40: was 2696
$ws.Cells.Item(40, 11).Value = 2261.6667  # K40: was 2696
$ws.Cells.Item(40, 13).Value = -2086.6667  # M40: was -2521

$ws.Cells.Item(64, 8).Value = 5047.25  # H64: was 5047.5
$ws.Cells.Item(64, 9).Value = 4994  # I64: was 0
$ws.Cells.Item(64, 10).Value = 5065  # J64: was 5047.5
$ws.Cells.Item(64, 11).Value = 4994  # K64: was 0
$ws.Cells.Item(64, 12).Value = 5065  # L64: was 5047.5
$ws.Cells.Item(64, 13).Value = -4746  # M64: was None
$ws.Cells.Item(64, 14).Value = -5561  # N64: was -5543.5

$ws.Cells.Item(67, 8).Value = 5047.25  # H67: was 5047.5
$ws.Cells.Item(67, 9).Value = 4994  # I67: was 0
$ws.Cells.Item(67, 10).Value = 5065  # J67: was 5047.5
$ws.Cells.Item(67, 11).Value = 4994  # K67: was 0
$ws.Cells.Item(67, 12).Value = 5065  # L67: was 5047.5
$ws.Cells.Item(67, 13).Value = -4136  # M67: was None
$ws.Cells.Item(67, 14).Value = -6781  # N67: was -6763.5

$ws.Cells.Item(106, 8).Value = 11478.667  # H106: was 9837.429
$ws.Cells.Item(106, 9).Value = 8968.25  # I106: was 6477
$ws.Cells.Item(106, 10).Value = 16499.5  # J106: was 30000
$ws.Cells.Item(106, 11).Value = 8968.25  # K106: was 6477
$ws.Cells.Item(106, 12).Value = 16499.5  # L106: was 30000
$ws.Cells.Item(106, 13).Value = -8337.25  # M106: was -5846
$ws.Cells.Item(106, 14).Value = -17761.5  # N106: was -31262

$ws.Cells.Item(127, 8).Value = 2879.6  # H127: was 2378.2856
$ws.Cells.Item(127, 9).Value = 2774.5  # I127: was 2224.6667
$ws.Cells.Item(127, 11).Value = 8323.5  # K127: was 6674.000100000001
$ws.Cells.Item(127, 13).Value = -3363.5  # M127: was -1714.000100000001

$ws.Cells.Item(129, 8).Value = 675.51514  # H129: was 684.5454999999999
$ws.Cells.Item(129, 9).Value = 945.53845  # I129: was 968.46155
$ws.Cells.Item(129, 11).Value = 2836.61535  # K129: was 2905.38465
$ws.Cells.Item(129, 13).Value = 2163.38465  # M129: was 2094.61535

$ws.Cells.Item(131, 8).Value = 1185.4  # H131: was 1142.25
$ws.Cells.Item(131, 9).Value = 1185.4  # I131: was 1142.25
$ws.Cells.Item(131, 11).Value = 3556.2  # K131: was 3426.75
$ws.Cells.Item(131, 13).Value = 1483.8  # M131: was 1613.25

$ws.Cells.Item(137, 8).Value = 1790.7333  # H137: was 1601.5834
$ws.Cells.Item(137, 9).Value = 1052.1818  # I137: was 1133.2
$ws.Cells.Item(137, 10).Value = 3821.75  # J137: was 3943.5
$ws.Cells.Item(137, 11).Value = 3156.5454  # K137: was 3399.6
$ws.Cells.Item(137, 12).Value = 11465.25  # L137: was 11830.5
$ws.Cells.Item(137, 13).Value = -606.5454  # M137: was -849.6000000000004
$ws.Cells.Item(137, 14).Value = -16565.25  # N137: was -16930.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 39727.535  # H32: was 49272.125
$ws.Cells.Item(32, 9).Value = 48035.914  # I32: was 64443
$ws.Cells.Item(32, 11).Value = 48035.914  # K32: was 64443
$ws.Cells.Item(32, 13).Value = -47748.914  # M32: was -64156

$ws.Cells.Item(61, 8).Value = 3822.4783  # H61: was 3898.087
$ws.Cells.Item(61, 9).Value = 3822.4783  # I61: was 3898.087
$ws.Cells.Item(61, 11).Value = 3822.4783  # K61: was 3898.087
$ws.Cells.Item(61, 13).Value = -3610.4783  # M61: was -3686.087

$ws.Cells.Item(97, 8).Value = 4688.457  # H97: was 4805.7646
$ws.Cells.Item(97, 9).Value = 5547.826  # I97: was 5765.909
$ws.Cells.Item(97, 10).Value = 3041.3333  # J97: was 3045.5
$ws.Cells.Item(97, 11).Value = 5547.826  # K97: was 5765.909
$ws.Cells.Item(97, 12).Value = 3041.3333  # L97: was 3045.5
$ws.Cells.Item(97, 13).Value = -5051.826  # M97: was -5269.909
$ws.Cells.Item(97, 14).Value = -4033.3333  # N97: was -4037.5

$ws.Cells.Item(136, 8).Value = 3822.4783  # H136: was 3898.087
$ws.Cells.Item(136, 9).Value = 3822.4783  # I136: was 3898.087
$ws.Cells.Item(136, 11).Value = 11467.4349  # K136: was 11694.261
$ws.Cells.Item(136, 13).Value = -8917.4349  # M136: was -9144.261

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1695.6471  # H86: was 1816.3667
$ws.Cells.Item(86, 9).Value = 1417.2727  # I86: was 1522.0526
$ws.Cells.Item(86, 10).Value = 2206  # J86: was 2324.7273
$ws.Cells.Item(86, 11).Value = 1417.2727  # K86: was 1522.0526
$ws.Cells.Item(86, 12).Value = 2206  # L86: was 2324.7273
$ws.Cells.Item(86, 13).Value = -294.2727  # M86: was -399.0526
$ws.Cells.Item(86, 14).Value = -4452  # N86: was -4570.7273

$ws.Cells.Item(89, 8).Value = 1695.6471  # H89: was 1816.3667
$ws.Cells.Item(89, 9).Value = 1417.2727  # I89: was 1522.0526
$ws.Cells.Item(89, 10).Value = 2206  # J89: was 2324.7273
$ws.Cells.Item(89, 11).Value = 7086.363499999999  # K89: was 7610.263
$ws.Cells.Item(89, 12).Value = 11030  # L89: was 11623.6365
$ws.Cells.Item(89, 13).Value = -1470.363499999999  # M89: was -1994.263
$ws.Cells.Item(89, 14).Value = -22262  # N89: was -22855.6365

$ws.Cells.Item(94, 8).Value = 1265.6522  # H94: was 1265.7826
$ws.Cells.Item(94, 9).Value = 1127.95  # I94: was 1128.1
$ws.Cells.Item(94, 11).Value = 1127.95  # K94: was 1128.1
$ws.Cells.Item(94, 13).Value = -676.95  # M94: was -677.0999999999999

$ws.Cells.Item(107, 8).Value = 1499.3334  # H107: was 1779.2858
$ws.Cells.Item(107, 9).Value = 666.3333  # I107: was 891.6
$ws.Cells.Item(107, 10).Value = 3998.3333  # J107: was 3998.5
$ws.Cells.Item(107, 11).Value = 666.3333  # K107: was 891.6
$ws.Cells.Item(107, 12).Value = 3998.3333  # L107: was 3998.5
$ws.Cells.Item(107, 13).Value = 1253.6667  # M107: was 1028.4
$ws.Cells.Item(107, 14).Value = -7838.3333  # N107: was -7838.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 55076.26  # H58: was 55128
$ws.Cells.Item(58, 9).Value = 61173.47  # I58: was 61231.293
$ws.Cells.Item(58, 11).Value = 61173.47  # K58: was 61231.293
$ws.Cells.Item(58, 13).Value = -60970.47  # M58: was -61028.293

$ws.Cells.Item(135, 8).Value = 85000  # H135: was 0
$ws.Cells.Item(135, 10).Value = 85000  # J135: was 0
$ws.Cells.Item(135, 12).Value = 85000  # L135: was 0
$ws.Cells.Item(135, 14).Value = -95140  # N135: was None

$ws.Cells.Item(136, 8).Value = 55076.26  # H136: was 55128
$ws.Cells.Item(136, 9).Value = 61173.47  # I136: was 61231.293
$ws.Cells.Item(136, 11).Value = 183520.41  # K136: was 183693.879
$ws.Cells.Item(136, 13).Value = -180970.41  # M136: was -181143.879

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 61988.89  # H37: was 62250
$ws.Cells.Item(37, 10).Value = 61988.89  # J37: was 62250
$ws.Cells.Item(37, 12).Value = 185966.67  # L37: was 186750
$ws.Cells.Item(37, 14).Value = -186190.67  # N37: was -186974

$ws.Cells.Item(50, 8).Value = 83547.836  # H50: was 100231
$ws.Cells.Item(50, 9).Value = 137.25  # I50: was 142.5
$ws.Cells.Item(50, 11).Value = 411.75  # K50: was 427.5
$ws.Cells.Item(50, 13).Value = 69.25  # M50: was 53.5

$ws.Cells.Item(53, 8).Value = 83547.836  # H53: was 100231
$ws.Cells.Item(53, 9).Value = 137.25  # I53: was 142.5
$ws.Cells.Item(53, 11).Value = 411.75  # K53: was 427.5
$ws.Cells.Item(53, 13).Value = 69.25  # M53: was 53.5

$ws.Cells.Item(58, 8).Value = 5568.6665  # H58: was 5733.3335
$ws.Cells.Item(58, 10).Value = 7752.5  # J58: was 7999.5
$ws.Cells.Item(58, 12).Value = 23257.5  # L58: was 23998.5
$ws.Cells.Item(58, 14).Value = -23513.5  # N58: was -24254.5

$ws.Cells.Item(86, 8).Value = 891.6667  # H86: was 713.1667
$ws.Cells.Item(86, 9).Value = 875  # I86: was 592
$ws.Cells.Item(86, 10).Value = 900  # J86: was 834.3333
$ws.Cells.Item(86, 11).Value = 2625  # K86: was 1776
$ws.Cells.Item(86, 12).Value = 2700  # L86: was 2502.9999
$ws.Cells.Item(86, 13).Value = -1439  # M86: was -590
$ws.Cells.Item(86, 14).Value = -5072  # N86: was -4874.9999

$ws.Cells.Item(89, 8).Value = 891.6667  # H89: was 713.1667
$ws.Cells.Item(89, 9).Value = 875  # I89: was 592
$ws.Cells.Item(89, 10).Value = 900  # J89: was 834.3333
$ws.Cells.Item(89, 11).Value = 7875  # K89: was 5328
$ws.Cells.Item(89, 12).Value = 8100  # L89: was 7508.9997
$ws.Cells.Item(89, 13).Value = -1947  # M89: was 600
$ws.Cells.Item(89, 14).Value = -19956  # N89: was -19364.9997

$ws.Cells.Item(140, 8).Value = 2703  # H140: was 2924.75
$ws.Cells.Item(140, 9).Value = 2703  # I140: was 2924.75
$ws.Cells.Item(140, 11).Value = 8109  # K140: was 8774.25
$ws.Cells.Item(140, 13).Value = -2929  # M140: was -3594.25

$ws.Cells.Item(141, 8).Value = 1250  # H141: was 999.6667
$ws.Cells.Item(141, 9).Value = 1250  # I141: was 999.6667
$ws.Cells.Item(141, 11).Value = 3750  # K141: was 2999.0001
$ws.Cells.Item(141, 13).Value = 1430  # M141: was 2180.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1988  # H102: was 2095.375
$ws.Cells.Item(102, 9).Value = 1186.4667  # I102: was 1251.9286
$ws.Cells.Item(102, 11).Value = 1186.4667  # K102: was 1251.9286
$ws.Cells.Item(102, 13).Value = 435.5333000000001  # M102: was 370.0714

$ws.Cells.Item(122, 8).Value = 1929.7812  # H122: was 1959.7742
$ws.Cells.Item(122, 10).Value = 3401.2  # J122: was 4001.5
$ws.Cells.Item(122, 12).Value = 10203.6  # L122: was 12004.5
$ws.Cells.Item(122, 14).Value = -15103.6  # N122: was -16904.5

$ws.Cells.Item(126, 8).Value = 9091.4  # H126: was 9919.357
$ws.Cells.Item(126, 9).Value = 7330.222  # I126: was 8559
$ws.Cells.Item(126, 11).Value = 21990.666  # K126: was 25677
$ws.Cells.Item(126, 13).Value = -19520.666  # M126: was -23207

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 8376.261  # H46: was 8902.951999999999
$ws.Cells.Item(46, 9).Value = 15054.143  # I46: was 17096.5
$ws.Cells.Item(46, 10).Value = 5454.6875  # J46: was 5625.533
$ws.Cells.Item(46, 11).Value = 15054.143  # K46: was 17096.5
$ws.Cells.Item(46, 12).Value = 5454.6875  # L46: was 5625.533
$ws.Cells.Item(46, 13).Value = -14866.143  # M46: was -16908.5
$ws.Cells.Item(46, 14).Value = -5830.6875  # N46: was -6001.533

$ws.Cells.Item(60, 8).Value = 0  # H60: was 54961
$ws.Cells.Item(60, 10).Value = 0  # J60: was 54961
$ws.Cells.Item(60, 12).Value = 0  # L60: was 54961
$ws.Cells.Item(60, 14).ClearContents()  # N60: was -55979

$ws.Cells.Item(68, 8).Value = 3297.5  # H68: was 4298
$ws.Cells.Item(68, 9).Value = 1800  # I68: was 1995
$ws.Cells.Item(68, 10).Value = 3597  # J68: was 4873.75
$ws.Cells.Item(68, 11).Value = 1800  # K68: was 1995
$ws.Cells.Item(68, 12).Value = 3597  # L68: was 4873.75
$ws.Cells.Item(68, 13).Value = -1051  # M68: was -1246
$ws.Cells.Item(68, 14).Value = -5095  # N68: was -6371.75

$ws.Cells.Item(71, 8).Value = 3297.5  # H71: was 4298
$ws.Cells.Item(71, 9).Value = 1800  # I71: was 1995
$ws.Cells.Item(71, 10).Value = 3597  # J71: was 4873.75
$ws.Cells.Item(71, 11).Value = 9000  # K71: was 9975
$ws.Cells.Item(71, 12).Value = 17985  # L71: was 24368.75
$ws.Cells.Item(71, 13).Value = -5256  # M71: was -6231
$ws.Cells.Item(71, 14).Value = -25473  # N71: was -31856.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 75817  # H62: was 75490.17999999999
$ws.Cells.Item(62, 9).Value = 9205.714  # I62: was 8610.5
$ws.Cells.Item(62, 10).Value = 122444.9  # J62: was 134938.78
$ws.Cells.Item(62, 11).Value = 9205.714  # K62: was 8610.5
$ws.Cells.Item(62, 12).Value = 122444.9  # L62: was 134938.78
$ws.Cells.Item(62, 13).Value = -8581.714  # M62: was -7986.5
$ws.Cells.Item(62, 14).Value = -123692.9  # N62: was -136186.78

$ws.Cells.Item(65, 8).Value = 75817  # H65: was 75490.17999999999
$ws.Cells.Item(65, 9).Value = 9205.714  # I65: was 8610.5
$ws.Cells.Item(65, 10).Value = 122444.9  # J65: was 134938.78
$ws.Cells.Item(65, 11).Value = 46028.57  # K65: was 43052.5
$ws.Cells.Item(65, 12).Value = 612224.5  # L65: was 674693.9
$ws.Cells.Item(65, 13).Value = -42908.57  # M65: was -39932.5
$ws.Cells.Item(65, 14).Value = -618464.5  # N65: was -680933.9

$ws.Cells.Item(81, 8).Value = 7405.4443  # H81: was 7054.9
$ws.Cells.Item(81, 10).Value = 9615  # J81: was 8798.571
$ws.Cells.Item(81, 12).Value = 19230  # L81: was 17597.142
$ws.Cells.Item(81, 14).Value = -21352  # N81: was -19719.142

$ws.Cells.Item(84, 8).Value = 7405.4443  # H84: was 7054.9
$ws.Cells.Item(84, 10).Value = 9615  # J84: was 8798.571
$ws.Cells.Item(84, 12).Value = 96150  # L84: was 87985.70999999999
$ws.Cells.Item(84, 14).Value = -106758  # N84: was -98593.70999999999
